# Auto-generated Excel COM-interop edit script
# Updates the cryptos price/volume table to the latest snapshot values
# (commit: "Updated cryptos list on Fri Oct 13 06:48:47 UTC 2023 with GitHub Actions")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.934.70'
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.546.19'
$ws.Range("E3").Value = '  -1.27%  '
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '205.70'
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.484'
$ws.Range("E6").Value = '  -0.34%  '
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '21.40'
$ws.Range("E9").Value = '  -2.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0582'
$ws.Range("E10").Value = '  -0.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0857'
$ws.Range("E11").Value = '  -0.66%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.765.69'
$ws.Range("E12").Value = '  -1.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.545.53'
$ws.Range("E13").Value = '  -1.28%  '
$ws.Range("E14").Value = '  -1.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.511'
$ws.Range("E15").Value = '  -0.83%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.902.09'
$ws.Range("E16").Value = '  -0.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.61'
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '213.79'
$ws.Range("E18").Value = '  -0.87%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0684'
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.20'
$ws.Range("E20").Value = '  -2.34%  '
$ws.Range("E21").Value = '  +0.35%  '
$ws.Range("E22").Value = '  -3.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.17'
$ws.Range("E23").Value = '  -0.48%  '
$ws.Range("E24").Value = '  -3.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.31'
$ws.Range("E25").Value = '  -0.46%  '
$ws.Range("E26").Value = '  -1.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.82'
$ws.Range("E27").Value = '  -0.95%  '
$ws.Range("E28").Value = '  +0.30%  '
$ws.Range("E29").Value = '  -0.33%  '
$ws.Range("E30").Value = '  -1.67%  '
$ws.Range("E31").Value = '  -1.32%  '
$ws.Range("E32").Value = '  +1.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.365.97'
$ws.Range("E33").Value = '  -2.34%  '
$ws.Range("E34").Value = '  +0.59%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.54'
$ws.Range("E35").Value = '  +0.42%  '
$ws.Range("E36").Value = '  +5.69%  '
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("E39").Value = '  -2.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.805'
$ws.Range("E40").Value = '  -1.52%  '
$ws.Range("E41").Value = '  +0.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.988'
$ws.Range("E42").Value = '  -0.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.48'
$ws.Range("E43").Value = '  -1.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.22'
$ws.Range("E44").Value = '  +1.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.41'
$ws.Range("E45").Value = '  -0.78%  '
$ws.Range("E46").Value = '  -2.88%  '
$ws.Range("B47").Value = 'mCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.25'
$ws.Range("E47").Value = '  -2.64%  '
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.679.24'
$ws.Range("E48").Value = '  -1.44%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.31'
$ws.Range("E49").Value = '  -0.68%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0506'
$ws.Range("E50").Value = '  +0.89%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₇0965'
$ws.Range("E51").Value = '  -1.54%  '
